$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old student-roster data (A1:G5) before laying down the new
# grade-book content, since the new layout leaves several cells blank where
# the old sheet had values (e.g. D2, H2, I2, B3, B5, E5, etc).
$ws.Range("A1:I6").ClearContents()

# Header row
$ws.Range("A1").Value = "Họ tên"
$ws.Range("B1").Value = "Điểm miệng cột 1"
$ws.Range("C1").Value = "Điểm miệng cột 2"
$ws.Range("D1").Value = "Điểm 15p cột 1"
$ws.Range("E1").Value = "Điểm 15p cột 2"
$ws.Range("F1").Value = "Điểm 1 tiết cột 1"
$ws.Range("G1").Value = "Điểm 1 tiết cột 2"
$ws.Range("H1").Value = "Điểm Giữa kỳ"
$ws.Range("I1").Value = "Điểm Cuối kỳ"

# Row 2 - Lý Anh Hiển
$ws.Range("A2").Value = "Lý Anh Hiển"
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 1
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1

# Row 3 - Nguyễn Dương Thanh Trúc
$ws.Range("A3").Value = "Nguyễn Dương Thanh Trúc"
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 2
$ws.Range("G3").Value = 2
$ws.Range("H3").Value = 2

# Row 4 - Nguyễn Văn G
$ws.Range("A4").Value = "Nguyễn Văn G"
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 1

# Row 5 - Lý Suni
$ws.Range("A5").Value = "Lý Suni"
$ws.Range("C5").Value = 3
$ws.Range("D5").Value = 4
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 3

# Row 6 - Lý Money
$ws.Range("A6").Value = "Lý Money"
